# Add a new "2022-Q3" worksheet right after "总计" (the first sheet),
# populate it with the quarterly fund-holding data, and update the
# "总计" (summary) sheet with a new leading row for 2022-Q3, shifting
# the existing quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" sheet right after "总计".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Header row.
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q3.Cells.Item(1, $i + 2)
    $cell.Value = $headers[$i]
    $cell.Style = $total.Cells.Item(1, 2).Style
}

# Data rows: index, fund code, fund name, fund scale, stock position total,
# position ratio, market value held (100M), position rank.
$rows = @(
    @(0,  "012930", "中庚价值先锋股票",                  "47.83", "94.71", "4.77", "2.2815", 5),
    @(1,  "007130", "中庚小盘价值股票",                  "75.87", "93.06", "2.57", "1.9499", 9),
    @(2,  "001955", "中欧养老产业混合A",                 "22.51", "92.70", "7.13", "1.6050", 9),
    @(3,  "010429", "中欧睿见混合A",                     "18.74", "91.22", "6.84", "1.2818", 9),
    @(4,  "011710", "中欧睿泽混合A",                     "8.33",  "90.88", "7.72", "0.6431", 9),
    @(5,  "004616", "中欧电子信息产业沪港深股票A",       "5.01",  "92.97", "5.78", "0.2896", 4),
    @(6,  "005763", "中欧电子信息产业沪港深股票C",       "3.88",  "92.97", "5.78", "0.2243", 4),
    @(7,  "012778", "中欧养老产业混合C",                 "2.80",  "92.70", "7.13", "0.1996", 9),
    @(8,  "002450", "平安睿享文娱灵活配置混合A",         "3.64",  "88.33", "3.67", "0.1336", 8),
    @(9,  "011711", "中欧睿泽混合C",                     "0.86",  "90.88", "7.72", "0.0664", 9),
    @(10, "002451", "平安睿享文娱灵活配置混合C",         "1.72",  "88.33", "3.67", "0.0631", 8),
    @(11, "011093", "永赢宏泽一年定期开放灵活配置混合",  "14.85", "47.39", "0.37", "0.0549", 8),
    @(12, "004351", "汇丰晋信珠三角区域发展混合",        "0.42",  "93.94", "3.37", "0.0142", 7),
    @(13, "015481", "中欧睿见混合C",                     "0.08",  "91.22", "6.84", "0.0055", 9)
)

foreach ($row in $rows) {
    $r = [int]$row[0] + 2

    $idxCell = $q3.Cells.Item($r, 1)
    $idxCell.Value = [int]$row[0]
    $idxCell.Style = $total.Cells.Item(2, 1).Style

    # Fund code is numeric-looking (leading zeros matter) -> force text.
    $codeCell = $q3.Cells.Item($r, 2)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row[1]

    # Fund name is plain text already.
    $q3.Cells.Item($r, 3).Value = $row[2]

    # Numeric-looking figures must remain text (leading zeros / fixed
    # decimals matter), so force a text number-format before assigning.
    for ($c = 4; $c -le 7; $c++) {
        $cell = $q3.Cells.Item($r, $c)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c - 1]
    }

    # Position rank is a genuine number.
    $q3.Cells.Item($r, 8).Value = [int]$row[7]
}

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: insert a new first data row for
#    2022-Q3 and push the existing quarters down by one row.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 1).Style = $total.Cells.Item(3, 1).Style
$total.Cells.Item(2, 2).Value = "2022-Q3"
$total.Cells.Item(2, 3).Value = 14
$total.Cells.Item(2, 4).Value = 8.81

# Renumber the index column (A) for the rows that shifted down.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
